$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AccountCheck")

# --- New columns V (Username2 bank account) and W (Password2) ---
# Write the repeating "AccountChek.bank18" / "bank18" values first so they
# take shared-string indices 77/78, then the header strings "Username2" /
# "Password2" so they land on 79/80 (matches the authored file's string order).
$ws.Range("V2:V6").Value = "AccountChek.bank18"
$ws.Range("W2:W6").Value = "bank18"
$ws.Range("V1").Value = "Username2"
$ws.Range("W1").Value = "Password2"

# --- Columns T (Username) / U (Password) also get filled in for rows 3-6 ---
# (rows 1 and 2 already had these values before this edit)
$ws.Range("T3:T6").Value = "Account.bank2"
$ws.Range("U3:U6").Value = "bank2"

# Column V was auto-fit by Excel to the widest entry ("AccountChek.bank18").
$ws.Columns.Item(22).ColumnWidth = 18.8

# Select the newly added block, matching the author's last selection.
$ws.Range("V2:W6").Select()
